$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    25.55153265035599,
    21.37877046517337,
    18.50481434201444,
    16.56452848725681,
    20.89329777001549,
    14.08859588408308,
    16.29753665093043,
    16.97428925808467,
    15.90832419965685,
    15.30353978893345,
    20.67429054105014,
    16.31714536707528,
    17.31691185864173,
    23.51107854801813,
    22.76608504505725,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062,
    25.42651250689062
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
